# Actualización automática hashcode mié feb  5 01:41:09 CET 2020
#
# This script updates the "hashcode" value column (column B) for a set of
# rows identified by their key in column A. For each entry below we verify
# that the row still holds the expected Key/Old value before writing the
# New value, and fall back to searching the whole sheet by key text if the
# row numbers ever shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 11; Key = "05-050301A"; Old = "287a03cb5ddcc8c51c66ebb6c60decca"; New = "10886ff983b31d2b23a61037bee03245" },
    @{ Row = 15; Key = "05-050207TP"; Old = "fef132f1f5ff59d4d500645fdae2eafc"; New = "995757dc0db9fd58a8f2be03086933e8" },
    @{ Row = 24; Key = "05-050316TC"; Old = "94b331c972a79e232e53636f1b848987"; New = "a555b1001cb4eb4774155e4e0d45ed42" },
    @{ Row = 29; Key = "05-050302A"; Old = "b3523c91a4da0c19819c5a321f6ac6e9"; New = "cee7288738f2d18a531208d0aff3184a" },
    @{ Row = 34; Key = "05-050316TP"; Old = "8136dbefd66d655fdc211f254fd091fe"; New = "199a671aa92b1cddc0ce99fa4e4b1e95" },
    @{ Row = 121; Key = "05-050301TP"; Old = "6f1993f7c3c6c8d0cc30ecd19cf4deca"; New = "4a8f025107cdceea31bb42aaa6359639" },
    @{ Row = 133; Key = "05-050312TP"; Old = "f978a7e8bfba99d00dc602838465a7b4"; New = "79d1d0c0b7ca3ccbf7b1e9c227371f3f" },
    @{ Row = 136; Key = "05-050312TC"; Old = "c84dd2a5765660273a59188852e315dd"; New = "ced153bef9faf7a242b0bc254c1cbd1b" },
    @{ Row = 159; Key = "05-050203TP"; Old = "86a32b40bf3869218dbb5318ac73dde7"; New = "3179b1019818ad8c556b64072a9463e9" },
    @{ Row = 162; Key = "05-050308A"; Old = "1abcbe56cddac5287b245556ff1850e4"; New = "fac192a900ed093137d7272371060418" },
    @{ Row = 169; Key = "05-050203TC"; Old = "75ad2a5365ea8a72ca5ddbbc28b828fb"; New = "d7bac45005bab7986cd39efac771ac50" },
    @{ Row = 175; Key = "05-050303TP"; Old = "5c1f0b2e58f0e164d0ce0420316e2c32"; New = "193cd4c5a8bab3fca87960a3d4334401" },
    @{ Row = 180; Key = "05-050303TC"; Old = "f490f3daff24fd8191d13c656d941609"; New = "5abe0996962ce49df8ad4ecad6d6e6b1" },
    @{ Row = 191; Key = "05-050314TP"; Old = "2a86ed6ba3d262609970061f0ef2d328"; New = "1083da5df02bf38f818a271508322574" },
    @{ Row = 198; Key = "05-050314TC"; Old = "95408bba2b4295e392962a87771f401e"; New = "8067240336eb47712eaecf0e3379c696" },
    @{ Row = 213; Key = "05-050303A"; Old = "84e1e733e8b824250a9bbc0e04afd984"; New = "3d3e8d23a97d243c3fb637cfccec89d7" },
    @{ Row = 339; Key = "05-050201TP"; Old = "885d675495acea9740f1c7bb31cfbbaa"; New = "addf5a747b264949fa9ae8e691ca5087" },
    @{ Row = 464; Key = "05-050204A"; Old = "c64fea71094245a6b65dbbf602a9480b"; New = "b189b4d6b4454b07494170016cc0a052" },
    @{ Row = 465; Key = "05-050313A"; Old = "538b5ea126904cb8272f36c92c2db2b9"; New = "3ac186dba6f835cc0dde39e9b9c8b581" },
    @{ Row = 477; Key = "03-030016A"; Old = "67996393d700dcfd73f87d83e57729d0"; New = "3f6233748c9d480d537076a8e25cd463" },
    @{ Row = 485; Key = "05-050314A"; Old = "6822ecb763f28c326f3b826b1a471ea6"; New = "3dab421690256830d891eb1dbd6545c8" },
    @{ Row = 507; Key = "05-050311A"; Old = "b2ed9656a757e542f1e2bbe43aef241e"; New = "1cacb3cee02312b2a93c65a2a344c7bf" },
    @{ Row = 508; Key = "05-050208TP"; Old = "f4ecf7d3761c99fd246bf4d08bdd9a00"; New = "f3c88963c669908676b5f56d7b21598e" },
    @{ Row = 513; Key = "05-050306TP"; Old = "ed3b78a07e857a3e0c3505d86522448d"; New = "57453290a028d0832d2d6a87aba3f3d1" },
    @{ Row = 521; Key = "05-050317TC"; Old = "c548e1027ffd749494f527c35b418364"; New = "7bec1385342fed9aa75716535350b327" },
    @{ Row = 532; Key = "05-050317TP"; Old = "231d36d68b4a94dae202778d18e76688"; New = "bd765d93499de8a428406c20c7de6700" },
    @{ Row = 555; Key = "05-050201A"; Old = "781565fc03d4b8852605f066d47696e9"; New = "fe7d557384f7dd0eb1f7e33fbc4ec243" },
    @{ Row = 624; Key = "05-050204TP"; Old = "98d74cdd8f1992c38d3de5c4f237d050"; New = "0690257d524fa65e2c39a24884c7519e" },
    @{ Row = 635; Key = "05-050204TC"; Old = "d91be6043d4519e7a2106349ed286d2a"; New = "b984c87dcf8554dba12699230be4fd78" },
    @{ Row = 637; Key = "05-050302TP"; Old = "8e9432c6f1a6d267311c0d455bc24788"; New = "4ffb3ea8d532b90ba41ae1b4caeab26c" },
    @{ Row = 657; Key = "05-050313TP"; Old = "a615830084fea3d6d3e77195e55486a6"; New = "260037e8ac135edabbd8ecebf727324e" },
    @{ Row = 663; Key = "05-050313TC"; Old = "2a7a243a38676356729aa6f06f136bb2"; New = "751c2a1d06e1e2780dafb73fa549ce0a" },
    @{ Row = 673; Key = "05-050208A"; Old = "003bd1a3349afac2db993828b457c703"; New = "e17400e9e4d237ac6e5cbbdb0c179f19" },
    @{ Row = 674; Key = "05-050317A"; Old = "9738a9d44cfe664e588c55837e54b311"; New = "5e8a3bfaf7d985c7619ed91006c40ba5" },
    @{ Row = 712; Key = "05-050315A"; Old = "9303bdfee0f12862e7f0b4458b662573"; New = "eed96e3ef1c25fb650d56cd4b8d8dc26" },
    @{ Row = 737; Key = "05-050316A"; Old = "c8e9cf10797a1c549d6b452f2ae90aef"; New = "7ac51dffd4b9e5f46303f624a41708a5" },
    @{ Row = 741; Key = "05-050207A"; Old = "4cd8d12abb0ad061a5045bafd15a0c72"; New = "4a3bb8dbe3e9fcb7011590db0761a9cd" },
    @{ Row = 750; Key = "05-050315TP"; Old = "5779b7978d3887cef77e946eaf833c0d"; New = "2958a31aa257329ad526cdbdc3e9be0e" },
    @{ Row = 827; Key = "05-050202TP"; Old = "7d618c8bf09746d171da3abda4a9112e"; New = "2c0679dd1930df1d151a7185272fd226" },
    @{ Row = 838; Key = "05-050311TC"; Old = "34263c36b3f531ee82374d4e0171ba28"; New = "609e13c97c8ea9422fcd925b50c0bb4f" },
    @{ Row = 843; Key = "05-050311TP"; Old = "34103d25a87755bfcc80a60448982bd0"; New = "2ee6460c61db675a0c438b7cc8ca8745" },
    @{ Row = 862; Key = "05-050309TC"; Old = "826863c9ee05d826fd66d3bbfcf566ca"; New = "2549441feec73fad726ef2286fad0e82" },
    @{ Row = 870; Key = "05-050309TP"; Old = "6280fb4a4df7ac6228b18992a8d15661"; New = "2ac976d0abbdb6753b1e8028cc220b23" },
    @{ Row = 963; Key = "03-030016TC"; Old = "74ce11a521c514d8df914174f6efb73d"; New = "3f574683856d8cc29639b08f7ab41e07" },
    @{ Row = 967; Key = "03-030016TP"; Old = "1382dc1aa6457e2dfe23d4db3af80247"; New = "ec8951b0c90004edf34c721157014b9d" }
)

$usedRows = @{}
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1

foreach ($u in $updates) {
    $row = $u.Row
    $keyCell = $ws.Cells.Item($row, 1)
    $valCell = $ws.Cells.Item($row, 2)

    if (($keyCell.Value2 -ne $u.Key) -or ($valCell.Value2 -ne $u.Old)) {
        # Row layout shifted: locate the row by the key text in column A.
        $found = $false
        for ($r = 1; $r -le $lastRow; $r++) {
            if ($usedRows.ContainsKey($r)) { continue }
            $c = $ws.Cells.Item($r, 1)
            if ($c.Value2 -eq $u.Key) {
                $row = $r
                $valCell = $ws.Cells.Item($row, 2)
                $found = $true
                break
            }
        }
        if (-not $found) {
            continue
        }
    }

    $valCell.Value2 = $u.New
    $usedRows[$row] = $true
}
